# Auto-generated PowerShell COM-interop script
# Applies the cryptos.xlsx price/volume refresh from the commit
# "Updated cryptos list on Sat Feb 17 02:37:02 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.923.83"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.786.87"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.85"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.44"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.588"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.43"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.53"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "3.226.35"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "2.772.81"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.935"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "51.912.25"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").Value = "0.0₃0972"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.27"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.49"
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.162"
$ws.Range("E28").Value = "  +15.62%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.22"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.90"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.69"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0839"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.20"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.94"
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.20"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.56"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.43"
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.74"
$ws.Range("E45").Value = "  -9.85%  "
$ws.Range("D46").Value = "2.081.17"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.80"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").Value = "  -4.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("E51").Value = "  -2.07%  "
